$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '72.302.62'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.647.73'
$ws.Range('E3').Value = '  -1.05%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '590.75'
$ws.Range('E5').Value = '  -2.55%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '175.09'
$ws.Range('E6').Value = '  -2.83%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -1.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.645.84'
$ws.Range('E9').Value = '  -1.26%  '
$ws.Range('E10').Value = '  -3.25%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.169'
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.132.73'
$ws.Range('E14').Value = '  -0.58%  '
$ws.Range('E15').Value = '  -2.81%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '72.197.63'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.00'
$ws.Range('E17').Value = '  -2.59%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.704.54'
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.28'
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.99'
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '372.81'
$ws.Range('E21').Value = '  -2.70%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.16'
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.07'
$ws.Range('E23').Value = '  -0.55%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '71.26'
$ws.Range('E24').Value = '  -2.22%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.28'
$ws.Range('E26').Value = '  -3.14%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.71'
$ws.Range('E27').Value = '  -3.61%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.780.60'
$ws.Range('E28').Value = '  -0.76%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0₃0963'
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('E31').Value = '  -1.81%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '501.42'
$ws.Range('E32').Value = '  -4.87%  '
$ws.Range('E33').Value = '  -3.14%  '
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '161.59'
$ws.Range('E36').Value = '  -1.49%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '19.35'
$ws.Range('E37').Value = '  -1.01%  '
$ws.Range('E38').Value = '  +2.35%  '
$ws.Range('E39').Value = '  -1.09%  '
$ws.Range('E40').Value = '  -2.91%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E42').Value = '  -6.59%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.57'
$ws.Range('E43').Value = '  -2.78%  '
$ws.Range('E44').Value = '  -4.22%  '
$ws.Range('E45').Value = '  -2.47%  '
$ws.Range('E46').Value = '  -0.45%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '153.35'
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.551'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('E49').Value = '  -1.73%  '
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0750'
$ws.Range('E51').Value = '  -1.30%  '
